# Apply cryptocurrency price/ranking updates to the active worksheet.
# The "Price" column (D) holds numeric-looking values that must stay as
# literal text (e.g. "274.20", "0.001663", trailing zeros preserved), so
# each target cell is forced to text format before the value is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> hashtable of column letter -> new value.
$updates = @{
    2  = @{ D = "274.20" }
    3  = @{ D = "22.95" }
    4  = @{ D = "6.352" }
    6  = @{ D = "3.659" }
    7  = @{ D = "6.668" }
    8  = @{ D = "1.367" }
    10 = @{ D = "0.01377" }
    11 = @{ D = "0.1633" }
    12 = @{ D = "0.08333" }
    14 = @{ D = "0.03107" }

    15 = @{ B = "BitMartToken";            C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";         D = "0.09305";  E = "14BitMartTokenBMX" }
    16 = @{ B = "MCDex";                   C = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb";                    D = "3.881";    E = "15MCDexMCB" }
    17 = @{ B = "BitForexToken";           C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";         D = "0.001663"; E = "16BitForexTokenBF" }
    18 = @{ B = "CoinExToken";             C = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet";          D = "0.04778";  E = "17CoinExTokenCET" }
    19 = @{ B = "TigerCash";               C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                 D = "0.006342"; E = "18TigerCashTCH" }
    20 = @{ B = "HotbitToken";             C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";           D = "0.005555"; E = "19HotbitTokenHTBWorstin24h" }
    21 = @{ B = "BitKan";                  C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";               D = "0.001091"; E = "20BitKanKAN" }
    22 = @{ B = "NitroEx";                 C = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";                D = "0.0001501"; E = "21NitroExNTX" }
    23 = @{ B = "LEO";                     C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                   D = "3.725";    E = "22LEOLEO" }
    24 = @{ B = "BTSEToken";               C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";            D = "2.370";    E = "23BTSETokenBTSE" }
    25 = @{ B = "BitpandaEcosystemToken";  C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";   D = "0.3381";   E = "24BitpandaEcosystemTokenBEST" }
    26 = @{ B = "ProBitToken";             C = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob";              D = "0.1269";   E = "25ProBitTokenPROB" }

    40 = @{ D = "0.04688" }
    41 = @{ D = "0.007033" }
    42 = @{ D = "0.1163" }
    43 = @{ D = "0.003601"; E = "42CEJICEJI" }
    44 = @{ D = "0.01216" }
    45 = @{ D = "0.00006253" }
    47 = @{ D = "0.9003" }
    48 = @{ D = "0.03074" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $value = $cols[$col]
        if ($col -eq "D") {
            # Keep the numeric-looking price string as literal text (leading
            # apostrophe forces text entry), then restore the plain "Normal"
            # style so no stray text-number-format style sticks to the cell.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
